$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) TC_Subset_01 (currently the 2nd worksheet) - add UserName/Password
#    columns before the existing "Subset Name" column.
# ------------------------------------------------------------------
$tcSubset01 = $wb.Worksheets.Item(2)
$tcSubset01.Range("B1:C1").EntireColumn.Insert()

# Preserve the exact shared-string insertion order seen in the target
# workbook: AutomationTestSubset, UserName, Password, automationUser,
# unilog123##, userName.
$tcSubset01.Range("D2").Value = "AutomationTestSubset"
$tcSubset01.Range("B1").Value = "UserName"
$tcSubset01.Range("C1").Value = "Password"
$tcSubset01.Range("B2").Value = "automationUser"
$tcSubset01.Range("C2").Value = "unilog123##"

$tcSubset01.Range("B:B").ColumnWidth = 14.666666666666666
$tcSubset01.Range("C:C").ColumnWidth = 10.666666666666666

$tcSubset01.Range("A1:C2").Select()

# ------------------------------------------------------------------
# 2) Insert a brand-new worksheet "TC_Subset_02" right after
#    TC_Subset_01.
# ------------------------------------------------------------------
$tcSubset02 = $wb.Worksheets.Add($null, $tcSubset01)
$tcSubset02.Name = "TC_Subset_02"

$tcSubset02.Range("A1").Value = "Sl. No"
$tcSubset02.Range("B1").Value = "UserName"
$tcSubset02.Range("C1").Value = "Password"
$tcSubset02.Range("A2").Value = 1
$tcSubset02.Range("B2").Value = "automationUser"
$tcSubset02.Range("C2").Value = "unilog123##"

$tcSubset02.Range("B:B").ColumnWidth = 14.666666666666666
$tcSubset02.Range("C:C").ColumnWidth = 10.666666666666666

$tcSubset02.Range("B1").Select()

# ------------------------------------------------------------------
# 3) TC_Subset_08 (now the 4th worksheet after the insert above) -
#    same UserName/Password columns, but header text is "userName"
#    (lower case) here, and it becomes the active tab.
# ------------------------------------------------------------------
$tcSubset08 = $wb.Worksheets.Item(4)
$tcSubset08.Range("B1:C1").EntireColumn.Insert()

$tcSubset08.Range("D2").Value = "AutomationTestSubset"
$tcSubset08.Range("B1").Value = "userName"
$tcSubset08.Range("C1").Value = "Password"
$tcSubset08.Range("B2").Value = "automationUser"
$tcSubset08.Range("C2").Value = "unilog123##"

$tcSubset08.Range("B:B").ColumnWidth = 14.666666666666666
$tcSubset08.Range("C:C").ColumnWidth = 10.666666666666666

$tcSubset08.Range("D2").Select()
